# More Hunter perks icons.
# Fill in the "icon present?" column (D) for several Hunter perk rows on Sheet1,
# and move the active cell selection to D33.
#
# Shared-string order matters: "Check" must become index 118 and "As-is" must
# become index 119, so we set the cell that introduces "Check" (D43) before the
# cells that introduce "As-is" (D28/D37/D40). "OK" already exists (index 117)
# and is simply reused.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 first, so the new string "Check" is appended to the shared strings
# table before "As-is".
$ws.Range("D43").Value = "Check"

# Row 28: new cell, was blank.
$ws.Range("D28").Value = "As-is"

# Row 32: new cell, was blank.
$ws.Range("D32").Value = "OK"

# Row 37: existing "OK" changes to "As-is".
$ws.Range("D37").Value = "As-is"

# Row 39: new cell, was blank.
$ws.Range("D39").Value = "OK"

# Row 40: new cell, was blank.
$ws.Range("D40").Value = "As-is"

# Update the saved selection/active cell shown in the sheet view.
$ws.Range("D33").Select()
